# Natmi following Dr Hou advice
# Rebuild the Nodal -> Acvr2b LR-pairs table with the full 3x3 cross of
# Sending cluster x Target cluster (ECs, FAPs, sCs), recomputed stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{A="ECs";  B="Nodal"; C="Acvr2b"; D="ECs";  E=3; F=1;                  G=1.089247;           H=3.267741;           I=0.5826095452486192; J=0.5826095452486191; K=2; L=0.6666666666666666; M=1.539665666666667;  N=4.618997; O=0.3572088291809875;  P=0.3572088291809875;  Q=1.677076208419667;  R=15.093685875777;  S=0.2081132735279268;  T=0.2081132735279268},
    @{A="ECs";  B="Nodal"; C="Acvr2b"; D="FAPs"; E=3; F=1;                  G=1.089247;           H=3.267741;           I=0.5826095452486192; J=0.5826095452486191; K=3; L=1;                  M=1.452872333333333;  N=4.358617; O=0.3370724153789985;  P=0.3370724153789985;  Q=1.582536830466333;  R=14.242831474197;  S=0.196381606639812;   T=0.196381606639812},
    @{A="ECs";  B="Nodal"; C="Acvr2b"; D="sCs";  E=3; F=1;                  G=1.089247;           H=3.267741;           I=0.5826095452486192; J=0.5826095452486191; K=3; L=1;                  M=1.317729666666667;  N=3.953189; O=0.3057187554400141;  P=0.3057187554400141;  Q=1.435333086227667;  R=12.917997776049;  S=0.1781146650808804;  T=0.1781146650808804},
    @{A="FAPs"; B="Nodal"; C="Acvr2b"; D="ECs";  E=2; F=0.6666666666666666; G=0.5602416666666666; H=1.680725;           I=0.2996585188171233; J=0.2996585188171233; K=2; L=0.6666666666666666; M=1.539665666666667;  N=4.618997; O=0.3572088291809875;  P=0.3572088291809875;  Q=0.8625848592027778; R=7.763263732825001; S=0.1070406686607735;  T=0.1070406686607735},
    @{A="FAPs"; B="Nodal"; C="Acvr2b"; D="FAPs"; E=2; F=0.6666666666666666; G=0.5602416666666666; H=1.680725;           I=0.2996585188171233; J=0.2996585188171233; K=3; L=1;                  M=1.452872333333333;  N=4.358617; O=0.3370724153789985;  P=0.3370724153789985;  Q=0.8139596174805555; R=7.325636557325;    S=0.1010066207265808;  T=0.1010066207265808},
    @{A="FAPs"; B="Nodal"; C="Acvr2b"; D="sCs";  E=2; F=0.6666666666666666; G=0.5602416666666666; H=1.680725;           I=0.2996585188171233; J=0.2996585188171233; K=3; L=1;                  M=1.317729666666667;  N=3.953189; O=0.3057187554400141;  P=0.3057187554400141;  Q=0.7382470646694445; R=6.644223582025;    S=0.09161122942976897; T=0.09161122942976897},
    @{A="sCs";  B="Nodal"; C="Acvr2b"; D="ECs";  E=2; F=0.6666666666666666; G=0.2201116666666666; H=0.6603349999999999; I=0.1177319359342576; J=0.1177319359342576; K=2; L=0.6666666666666666; M=1.539665666666667;  N=4.618997; O=0.3572088291809875;  P=0.3572088291809875;  Q=0.3388983759994444; R=3.050085383995;    S=0.04205488699228718; T=0.04205488699228718},
    @{A="sCs";  B="Nodal"; C="Acvr2b"; D="FAPs"; E=2; F=0.6666666666666666; G=0.2201116666666666; H=0.6603349999999999; I=0.1177319359342576; J=0.1177319359342576; K=3; L=1;                  M=1.452872333333333;  N=4.358617; O=0.3370724153789985;  P=0.3370724153789985;  Q=0.3197941507438888; R=2.878147356695;    S=0.03968418801260571; T=0.03968418801260571},
    @{A="sCs";  B="Nodal"; C="Acvr2b"; D="sCs";  E=2; F=0.6666666666666666; G=0.2201116666666666; H=0.6603349999999999; I=0.1177319359342576; J=0.1177319359342576; K=3; L=1;                  M=1.317729666666667;  N=3.953189; O=0.3057187554400141;  P=0.3057187554400141;  Q=0.2900476731461111; R=2.610429058315;    S=0.0359928609293647;  T=0.0359928609293647}
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    foreach ($col in $cols) {
        $ws.Range($col + $r).Value = $row[$col]
    }
}

Write-Output "wrote $($data.Count) data rows"
